$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Remove the old column A (TAXON's numeric "11-texon id" column) entirely,
# shifting B:F left to A:E.
$ws.Range("A:A").Delete()

# Fix the mis-spaced header text ("MODEL_CONDITION" -> "MODELCONDITION"),
# now located in column D after the shift.
$ws.Range("D1").Value = "MODELCONDITION"
